$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.402.94'
$ws.Range('E2').Value = '  -4.78%  '
$ws.Range('D3').Value = '1.571.18'
$ws.Range('E3').Value = '  -4.81%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.ClearFormats()
$ws.Range('E4').Value = '  +0.31%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.24%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '292.02'
$c.ClearFormats()
$ws.Range('E6').Value = '  -2.74%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.3664'
$c.ClearFormats()
$ws.Range('E7').Value = '  -3.30%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '49.46'
$c.ClearFormats()
$ws.Range('E8').Value = '  -2.56%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.3337'
$c.ClearFormats()
$ws.Range('E9').Value = '  -6.70%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '1.160'
$c.ClearFormats()
$ws.Range('E10').Value = '  -5.34%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.07543'
$c.ClearFormats()
$ws.Range('E11').Value = '  -7.01%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '1.004'
$c.ClearFormats()
$ws.Range('E12').Value = '  +0.25%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '20.98'
$c.ClearFormats()
$ws.Range('E13').Value = '  -5.09%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '6.109'
$c.ClearFormats()
$ws.Range('E14').Value = '  -4.77%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '6.843'
$c.ClearFormats()
$ws.Range('E15').Value = '  -7.75%  '
$ws.Range('D16').Value = '1.568.87'
$ws.Range('E16').Value = '  -5.42%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.00001131'
$c.ClearFormats()
$ws.Range('E17').Value = '  -5.98%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '89.44'
$c.ClearFormats()
$ws.Range('E18').Value = '  -7.83%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.06754'
$c.ClearFormats()
$ws.Range('E19').Value = '  -3.45%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.ClearFormats()
$ws.Range('E20').Value = '  +0.26%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '6.225'
$c.ClearFormats()
$ws.Range('E21').Value = '  -8.29%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '16.26'
$c.ClearFormats()
$ws.Range('E22').Value = '  -6.97%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '11.87'
$c.ClearFormats()
$ws.Range('E23').Value = '  -5.87%  '
$ws.Range('D24').Value = '22.440.18'
$ws.Range('E24').Value = '  -4.69%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.384'
$c.ClearFormats()
$ws.Range('E25').Value = '  -4.54%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.893'
$c.ClearFormats()
$ws.Range('E26').Value = '  -1.33%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '19.68'
$c.ClearFormats()
$ws.Range('E27').Value = '  -6.24%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '145.61'
$c.ClearFormats()
$ws.Range('E28').Value = '  -4.20%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '4.943'
$c.ClearFormats()
$ws.Range('E29').Value = '  -5.61%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '124.68'
$c.ClearFormats()
$ws.Range('E30').Value = '  -6.30%  '
$ws.Range('D31').Value = '1.744.02'
$ws.Range('E31').Value = '  -5.10%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '6.204'
$c.ClearFormats()
$ws.Range('E32').Value = '  -11.34%  '
$ws.Range('B33').Value = 'WEMIXTOKEN'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '2.013'
$c.ClearFormats()
$ws.Range('E33').Value = '  -6.27%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.9583'
$c.ClearFormats()
$ws.Range('E34').Value = '  -8.07%  '
$ws.Range('B35').Value = 'FraxShare'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '10.28'
$c.ClearFormats()
$ws.Range('E35').Value = '  -13.64%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.08523'
$c.ClearFormats()
$ws.Range('E36').Value = '  -2.38%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.02502'
$c.ClearFormats()
$ws.Range('E37').Value = '  -8.72%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.2286'
$c.ClearFormats()
$ws.Range('E38').Value = '  -6.87%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.06542'
$c.ClearFormats()
$ws.Range('E39').Value = '  -4.97%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '5.438'
$c.ClearFormats()
$ws.Range('E40').Value = '  -9.25%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.260'
$c.ClearFormats()
$ws.Range('E41').Value = '  -4.41%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '11.72'
$c.ClearFormats()
$ws.Range('E42').Value = '  -11.01%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.6340'
$c.ClearFormats()
$ws.Range('E43').Value = '  -8.52%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '14.51'
$c.ClearFormats()
$ws.Range('E44').Value = '  -8.21%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.ClearFormats()
$ws.Range('E45').Value = '  +0.21%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '3.776'
$c.ClearFormats()
$ws.Range('E46').Value = '  -3.96%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.5960'
$c.ClearFormats()
$ws.Range('E47').Value = '  -7.70%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.102'
$c.ClearFormats()
$ws.Range('E48').Value = '  -7.34%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '121.55'
$c.ClearFormats()
$ws.Range('E49').Value = '  -4.50%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.07254'
$c.ClearFormats()
$ws.Range('E50').Value = '  -7.59%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.170'
$c.ClearFormats()
$ws.Range('E51').Value = '  -1.95%  '
